# Applies the "Updated cryptos list" data refresh (Mon Aug 28 06:30:06 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.101.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.69%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.651.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.79%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.41%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5207"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.006"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2618"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.88%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06284"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07797"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.470"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.95%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.646.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.878.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5535"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7997"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.099.23"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.625"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "194.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.940"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.007"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1203"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.157"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.49%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05697"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.266"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.480"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.359"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.591"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.16%  "

$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.798"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.98%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9491"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.412"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5655"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01587"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.950"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.056.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.006"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8417"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.790.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05368"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.79%  "

$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.53%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4398"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.41%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₈103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.941"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.94%  "
